# epexspot_prices.xlsx update
#  - "Prix Spot": a new daily column "08-nov" is inserted before the
#    "01-oct." block (i.e. before column DM), shifting the October block
#    one column to the right (DM:EQ -> DN:ER). The new column gets "-"
#    placeholders for every hourly row (no data yet for that day).
#  - "Gaz" / "CO2": a new trailing row is appended with the next day's
#    (2025-11-06) closing price.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": insert one day column before DM (new "08-nov" day)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Shift DM:EQ (and everything after) one column to the right, inserting a
# fresh blank column at DM.
$ws1.Range("DM1:DM25").EntireColumn.Insert()

# Header cell for the newly inserted column.
$ws1.Range("DM1").Value = "08-nov"

# No observations yet for 08-nov, so every hourly row gets the same "-"
# placeholder already used elsewhere in the sheet for missing data.
for ($r = 2; $r -le 25; $r++) {
    $ws1.Range("DM$r").Value = "-"
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append 2025-11-06 closing price
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Stash a text-producing formula in a scratch cell, copy it, and paste
# values-only onto the target cell. That keeps the date as literal text
# ("2025-11-06") instead of letting Excel's type inference reinterpret
# the typed string as a date serial number.
$scratch = $wsGaz.Range("D1")
$scratch.Formula = '="2025-11-06"'
$scratch.Copy() | Out-Null
$wsGaz.Range("A145").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null

$wsGaz.Range("B145").Value = 30.35

# ---------------------------------------------------------------------
# Sheet "CO2": append 2025-11-06 closing price
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$scratch2 = $wsCo2.Range("D1")
$scratch2.Formula = '="2025-11-06"'
$scratch2.Copy() | Out-Null
$wsCo2.Range("A145").PasteSpecial(-4163) | Out-Null
$scratch2.Clear() | Out-Null

$wsCo2.Range("B145").Value = 79.94
